$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "naiveAR2" / "ifoCast" column headers (D1 <-> E1 labels)
# and the corresponding error-column headers (G1 <-> H1 labels).
$ws.Range("D1").Value = "ifoCast"
$ws.Range("E1").Value = "naiveAR2"
$ws.Range("G1").Value = "error_realized_minus_ifoCast"
$ws.Range("H1").Value = "error_realized_minus_naiveAR2"

# Negate every value in the error columns F, G, H for data rows 2-47.
for ($r = 2; $r -le 47; $r++) {
    foreach ($col in @("F","G","H")) {
        $cell = $ws.Range("$col$r")
        $cell.Value = -1 * $cell.Value()
    }
}
